## R/data/quiz230925.xlsx -- append 16 new survey responses (rows 423-438)
## plus the associated sharedStrings / dimension / sheet-view bookkeeping that
## Excel updates automatically when new rows are entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Seed formatting for the new rows by cloning existing template rows ----
# Most rows answer Q7 in column M ("row 3" layout); a few answer in column N
# instead ("row 2" layout, which skips M). Clone whichever matches each row so
# the new cells pick up the right style indices without inventing new ones.
$ws.Range("A3:M3").Copy()
$ws.Range("A423:M438").PasteSpecial(-4122)

$ws.Range("A2:N2").Copy()
$ws.Range("A425:N425").PasteSpecial(-4122)
$ws.Range("A431:N431").PasteSpecial(-4122)
$ws.Range("A432:N432").PasteSpecial(-4122)
$ws.Range("A433:N433").PasteSpecial(-4122)
$ws.Range("A438:N438").PasteSpecial(-4122)

# The "row 2" template has no M column; drop the blank placeholder it leaves behind
# on these rows so they end up with an N cell only, matching the survey answer.
$ws.Cells.Item(425,13).Clear()
$ws.Cells.Item(431,13).Clear()
$ws.Cells.Item(432,13).Clear()
$ws.Cells.Item(433,13).Clear()
$ws.Cells.Item(438,13).Clear()

$ws.Application.CutCopyMode = $false
$ws.Range("A423:A438").EntireRow.RowHeight = 15.75

# ---- 2. Fill in the actual response data ----
$ws.Cells.Item(423,1).Value = 45200.93687467593
$ws.Cells.Item(423,2).Value = 'sumin102573@naver.com'
$ws.Cells.Item(423,3).Value = '경영학과'
$ws.Cells.Item(423,4).Value = 20212922
$ws.Cells.Item(423,5).Value = '김수민'
$ws.Cells.Item(423,6).Value = '민주 문자'
$ws.Cells.Item(423,7).Value = '한글'
$ws.Cells.Item(423,8).Value = '하나도 없다'
$ws.Cells.Item(423,9).Value = 0.9
$ws.Cells.Item(423,10).Value = '미국'
$ws.Cells.Item(423,11).Value = '건강이 좋지 않다'
$ws.Cells.Item(423,12).Value = 'Red'
$ws.Cells.Item(423,13).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(424,1).Value = 45200.93968752315
$ws.Cells.Item(424,2).Value = 'limmh96@gmail.com'
$ws.Cells.Item(424,3).Value = '광고홍보학과'
$ws.Cells.Item(424,4).Value = 20232632
$ws.Cells.Item(424,5).Value = '임민호'
$ws.Cells.Item(424,6).Value = '엘리트 문자'
$ws.Cells.Item(424,7).Value = '한글'
$ws.Cells.Item(424,8).Value = '2개'
$ws.Cells.Item(424,9).Value = 0.8
$ws.Cells.Item(424,10).Value = '미국'
$ws.Cells.Item(424,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(424,12).Value = 'Red'
$ws.Cells.Item(424,13).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(425,1).Value = 45200.94147167824
$ws.Cells.Item(425,2).Value = 'scw0922@naver.com'
$ws.Cells.Item(425,3).Value = '간호학과'
$ws.Cells.Item(425,4).Value = 20236256
$ws.Cells.Item(425,5).Value = '신채원'
$ws.Cells.Item(425,6).Value = '민주 문자'
$ws.Cells.Item(425,7).Value = '한글'
$ws.Cells.Item(425,8).Value = '하나도 없다'
$ws.Cells.Item(425,9).Value = 0.5
$ws.Cells.Item(425,10).Value = '미국'
$ws.Cells.Item(425,11).Value = '남들을 덜 신뢰한다'
$ws.Cells.Item(425,12).Value = 'Black'
$ws.Cells.Item(425,14).Value = '헐, 반 밖에 안 남았네.'

$ws.Cells.Item(426,1).Value = 45200.94627571759
$ws.Cells.Item(426,2).Value = 'harin3040@naver.com'
$ws.Cells.Item(426,3).Value = '심리학과'
$ws.Cells.Item(426,4).Value = 20232113
$ws.Cells.Item(426,5).Value = '김현진'
$ws.Cells.Item(426,6).Value = '엘리트 문자'
$ws.Cells.Item(426,7).Value = '한글'
$ws.Cells.Item(426,8).Value = '하나도 없다'
$ws.Cells.Item(426,9).Value = 0.8
$ws.Cells.Item(426,10).Value = '이탈리아'
$ws.Cells.Item(426,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(426,12).Value = 'Red'
$ws.Cells.Item(426,13).Value = '헐, 반 밖에 안 남았네.'

$ws.Cells.Item(427,1).Value = 45200.947467094906
$ws.Cells.Item(427,2).Value = 'shanesun0923@gmail.com'
$ws.Cells.Item(427,3).Value = '간호학과'
$ws.Cells.Item(427,4).Value = 20236253
$ws.Cells.Item(427,5).Value = '선세인'
$ws.Cells.Item(427,6).Value = '민주 문자'
$ws.Cells.Item(427,7).Value = '한글'
$ws.Cells.Item(427,8).Value = '3개'
$ws.Cells.Item(427,9).Value = 0.9
$ws.Cells.Item(427,10).Value = '대한민국'
$ws.Cells.Item(427,11).Value = '사회활동이나 자원활동에 덜 참여한다'
$ws.Cells.Item(427,12).Value = 'Red'
$ws.Cells.Item(427,13).Value = '모름/기타'

$ws.Cells.Item(428,1).Value = 45200.95396979166
$ws.Cells.Item(428,2).Value = 'sung93716@gmail.com'
$ws.Cells.Item(428,3).Value = '데이터사이언스학부'
$ws.Cells.Item(428,4).Value = 20233261
$ws.Cells.Item(428,5).Value = '한예림'
$ws.Cells.Item(428,6).Value = '민주 문자'
$ws.Cells.Item(428,7).Value = '한글'
$ws.Cells.Item(428,8).Value = '2개'
$ws.Cells.Item(428,9).Value = 0.1
$ws.Cells.Item(428,10).Value = '이탈리아'
$ws.Cells.Item(428,11).Value = '2배 정도 실직할 가능성이 높다'
$ws.Cells.Item(428,12).Value = 'Red'
$ws.Cells.Item(428,13).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(429,1).Value = 45200.95795295139
$ws.Cells.Item(429,2).Value = 'molly7624@naver.com'
$ws.Cells.Item(429,3).Value = '디지털미디어콘텐츠전공'
$ws.Cells.Item(429,4).Value = 20211516
$ws.Cells.Item(429,5).Value = '변재은'
$ws.Cells.Item(429,6).Value = '민주 문자'
$ws.Cells.Item(429,7).Value = '한글'
$ws.Cells.Item(429,8).Value = '1개'
$ws.Cells.Item(429,9).Value = 0.8
$ws.Cells.Item(429,10).Value = '대한민국'
$ws.Cells.Item(429,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(429,12).Value = 'Red'
$ws.Cells.Item(429,13).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(430,1).Value = 45200.96560776621
$ws.Cells.Item(430,2).Value = 'choe0119@gmail.com'
$ws.Cells.Item(430,3).Value = '의예과'
$ws.Cells.Item(430,4).Value = 20226176
$ws.Cells.Item(430,5).Value = '최태웅'
$ws.Cells.Item(430,6).Value = '엘리트 문자'
$ws.Cells.Item(430,7).Value = '한자'
$ws.Cells.Item(430,8).Value = '1개'
$ws.Cells.Item(430,9).Value = 0.2
$ws.Cells.Item(430,10).Value = '영국'
$ws.Cells.Item(430,11).Value = '건강이 좋지 않다'
$ws.Cells.Item(430,12).Value = 'Red'
$ws.Cells.Item(430,13).Value = '헐, 반 밖에 안 남았네.'

$ws.Cells.Item(431,1).Value = 45200.970568564815
$ws.Cells.Item(431,2).Value = 'dksdksqh1018@naver.com'
$ws.Cells.Item(431,3).Value = '미디어스쿨'
$ws.Cells.Item(431,4).Value = 20232549
$ws.Cells.Item(431,5).Value = '안보민'
$ws.Cells.Item(431,6).Value = '민주 문자'
$ws.Cells.Item(431,7).Value = '한글'
$ws.Cells.Item(431,8).Value = '하나도 없다'
$ws.Cells.Item(431,9).Value = 0.8
$ws.Cells.Item(431,10).Value = '대한민국'
$ws.Cells.Item(431,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(431,12).Value = 'Black'
$ws.Cells.Item(431,14).Value = '헐, 반 밖에 안 남았네.'

$ws.Cells.Item(432,1).Value = 45200.97358670139
$ws.Cells.Item(432,2).Value = 'gaejisub@gmail.com'
$ws.Cells.Item(432,3).Value = '콘텐츠it'
$ws.Cells.Item(432,4).Value = 20225169
$ws.Cells.Item(432,5).Value = '배승유'
$ws.Cells.Item(432,6).Value = '민주 문자'
$ws.Cells.Item(432,7).Value = '한글'
$ws.Cells.Item(432,8).Value = '하나도 없다'
$ws.Cells.Item(432,9).Value = 0.2
$ws.Cells.Item(432,10).Value = '미국'
$ws.Cells.Item(432,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(432,12).Value = 'Black'
$ws.Cells.Item(432,14).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(433,1).Value = 45200.98353868056
$ws.Cells.Item(433,2).Value = 'gustj1654@naver.com'
$ws.Cells.Item(433,3).Value = '심리학과'
$ws.Cells.Item(433,4).Value = 20232137
$ws.Cells.Item(433,5).Value = '조현서'
$ws.Cells.Item(433,6).Value = '민주 문자'
$ws.Cells.Item(433,7).Value = '한글'
$ws.Cells.Item(433,8).Value = '1개'
$ws.Cells.Item(433,9).Value = 0.8
$ws.Cells.Item(433,10).Value = '대한민국'
$ws.Cells.Item(433,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(433,12).Value = 'Black'
$ws.Cells.Item(433,14).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(434,1).Value = 45200.98789751157
$ws.Cells.Item(434,2).Value = 'yejin4259@naver.com'
$ws.Cells.Item(434,3).Value = '언어청각학부'
$ws.Cells.Item(434,4).Value = 20233951
$ws.Cells.Item(434,5).Value = '이예진'
$ws.Cells.Item(434,6).Value = '민주 문자'
$ws.Cells.Item(434,7).Value = '한글'
$ws.Cells.Item(434,8).Value = '1개'
$ws.Cells.Item(434,9).Value = 0.8
$ws.Cells.Item(434,10).Value = '대한민국'
$ws.Cells.Item(434,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(434,12).Value = 'Red'
$ws.Cells.Item(434,13).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(435,1).Value = 45200.99348971064
$ws.Cells.Item(435,2).Value = 'rhy0787@naver.com'
$ws.Cells.Item(435,3).Value = '식품영양학과'
$ws.Cells.Item(435,4).Value = 20213827
$ws.Cells.Item(435,5).Value = '유희영'
$ws.Cells.Item(435,6).Value = '민주 문자'
$ws.Cells.Item(435,7).Value = '한자'
$ws.Cells.Item(435,8).Value = '하나도 없다'
$ws.Cells.Item(435,9).Value = 0.8
$ws.Cells.Item(435,10).Value = '대한민국'
$ws.Cells.Item(435,11).Value = '시간당 중위 임금이 60% 낮다'
$ws.Cells.Item(435,12).Value = 'Red'
$ws.Cells.Item(435,13).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(436,1).Value = 45200.99731702547
$ws.Cells.Item(436,2).Value = 'kxjenlee@naver.com'
$ws.Cells.Item(436,3).Value = '글로벌비즈니스'
$ws.Cells.Item(436,4).Value = 20226417
$ws.Cells.Item(436,5).Value = '이제인'
$ws.Cells.Item(436,6).Value = '민주 문자'
$ws.Cells.Item(436,7).Value = '한글'
$ws.Cells.Item(436,8).Value = '1개'
$ws.Cells.Item(436,9).Value = 0.5
$ws.Cells.Item(436,10).Value = '영국'
$ws.Cells.Item(436,11).Value = '2배 정도 실직할 가능성이 높다'
$ws.Cells.Item(436,12).Value = 'Red'
$ws.Cells.Item(436,13).Value = '모름/기타'

$ws.Cells.Item(437,1).Value = 45201.00346453703
$ws.Cells.Item(437,2).Value = 'tjdus3641@gmail.com'
$ws.Cells.Item(437,3).Value = '간호학과'
$ws.Cells.Item(437,4).Value = 20226283
$ws.Cells.Item(437,5).Value = '장서연'
$ws.Cells.Item(437,6).Value = '민주 문자'
$ws.Cells.Item(437,7).Value = '한글'
$ws.Cells.Item(437,8).Value = '2개'
$ws.Cells.Item(437,9).Value = 0.8
$ws.Cells.Item(437,10).Value = '대한민국'
$ws.Cells.Item(437,11).Value = '사회활동이나 자원활동에 덜 참여한다'
$ws.Cells.Item(437,12).Value = 'Red'
$ws.Cells.Item(437,13).Value = '휴우, 그래도 반이나 남았네.'

$ws.Cells.Item(438,1).Value = 45201.00409956019
$ws.Cells.Item(438,2).Value = 'rkqls3333@gmail.com'
$ws.Cells.Item(438,3).Value = '간호학과'
$ws.Cells.Item(438,4).Value = 20236205
$ws.Cells.Item(438,5).Value = '권가빈'
$ws.Cells.Item(438,6).Value = '민주 문자'
$ws.Cells.Item(438,7).Value = '한글'
$ws.Cells.Item(438,8).Value = '2개'
$ws.Cells.Item(438,9).Value = 0.2
$ws.Cells.Item(438,10).Value = '대한민국'
$ws.Cells.Item(438,11).Value = '2배 정도 실직할 가능성이 높다'
$ws.Cells.Item(438,12).Value = 'Black'
$ws.Cells.Item(438,14).Value = '모름/기타'

# ---- 3. Restore the frozen-pane scroll position / active cell the author left ----
$ws.Application.Goto($ws.Range("A411"), $true)
$ws.Range("D444").Select()
